$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.004.65"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "3.447.84"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "3.443.77"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "4.044.29"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.08"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.81%  "
$ws.Range("D16").Value = "65.956.50"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "3.445.28"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.531"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").Value = "2.779.33"
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0686"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0291"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "328.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.26%  "
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.37%  "
